$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Task name change: "Daemon tasky" -> "Server tasky"
$ws.Range("B14").Value = "Server tasky"

# Status updates in column AE (task progress tracker)
$ws.Range("AE12").Value = "???"
$ws.Range("AE13").Value = "Odloženo"
$ws.Range("AE14").Value = "Dokončeno"
$ws.Range("AE17").Value = "???"
$ws.Range("AE18").Value = "Dokončeno"
$ws.Range("AE19").Value = "Dokončeno"
$ws.Range("AE20").Value = "Probíha"
$ws.Range("AE21").Value = "Probíha"
$ws.Range("AE22").Value = "Probíha"

# Update comment text on AE13 (first line changes, second line stays the same)
$comment = $ws.Range("AE13").Comment
$comment.Text("lukx:`nUdělám s Rambem.`nDruhy backapů byli přidány do DB")

# Update the active selection as last recorded in the file
$ws.Range("AG21").Select()
